$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 305
$ws.Cells.Item(305, 1).Value = 44796
$ws.Cells.Item(305, 2).Value = "KA04MM4818"
$ws.Cells.Item(305, 3).Value = "RITZ"
$ws.Cells.Item(305, 4).Value = "GENERAL CHECKUP         WW"
$ws.Cells.Item(305, 5).Value = "WORK IN PROGRESS"

# Row 306
$ws.Cells.Item(306, 1).Value = 44796
$ws.Cells.Item(306, 2).Value = "KA03MB5345"
$ws.Cells.Item(306, 3).Value = "WAGON R"
$ws.Cells.Item(306, 4).Value = "GENERAL CHECKUP"
$ws.Cells.Item(306, 5).Value = "WORK DONE "

# Row 307
$ws.Cells.Item(307, 1).Value = 44796
$ws.Cells.Item(307, 2).Value = "KA05NA0487"
$ws.Cells.Item(307, 3).Value = "TIAGO"
$ws.Cells.Item(307, 4).Value = "BODY SHOP"
$ws.Cells.Item(307, 5).Value = "WORK DONE"
$ws.Cells.Item(307, 6).Value = 31604
$ws.Cells.Item(307, 7).Value = "  INSURANCE"

# Row 308
$ws.Cells.Item(308, 1).Value = 44796
$ws.Cells.Item(308, 2).Value = "KA03MN9595"
$ws.Cells.Item(308, 3).Value = "FORTUNER"
$ws.Cells.Item(308, 4).Value = "GENERAL CHECKUP"
$ws.Cells.Item(308, 5).Value = "WORK DONE DELIVERED"
$ws.Cells.Item(308, 6).Value = 32292
$ws.Cells.Item(308, 7).Value = "CREDIT"

# Row 309
$ws.Cells.Item(309, 1).Value = 44796
$ws.Cells.Item(309, 2).Value = "KA51MF0652"
$ws.Cells.Item(309, 3).Value = "XUV 500"
$ws.Cells.Item(309, 4).Value = "BODY POLISHING"
$ws.Cells.Item(309, 5).Value = "WORK DONE DELIVERED"
$ws.Cells.Item(309, 6).Value = 6500
$ws.Cells.Item(309, 7).Value = "G PAY"

# Scroll / selection to match final view state
$excel.ActiveWindow.ScrollRow = 292
$excel.ActiveWindow.ScrollColumn = 1
$null = $ws.Range("H309").Select()
